$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (German)
$ws.Range("B2").Value = 6444
$ws.Range("C2").Value = 4452
$ws.Range("D2").Value = 5810
$ws.Range("E2").Value = 6189
$ws.Range("F2").Value = 6347

# Row 6 (Greek)
$ws.Range("B6").Value = 5936
$ws.Range("C6").Value = 2233
$ws.Range("D6").Value = 4158
$ws.Range("E6").Value = 5709
$ws.Range("F6").Value = 5930

# Row 11 (Finnish)
$ws.Range("B11").Value = 7926
$ws.Range("C11").Value = 5729
$ws.Range("D11").Value = 7921
$ws.Range("E11").Value = 7926
$ws.Range("F11").Value = 7926

# Row 12 (Basque)
$ws.Range("B12").Value = 4432
$ws.Range("C12").Value = 4313
$ws.Range("D12").Value = 4411
$ws.Range("E12").Value = 4419
$ws.Range("F12").Value = 4424

# Row 14 (Japanese)
$ws.Range("B14").Value = 36000
$ws.Range("C14").Value = 34046
$ws.Range("D14").Value = 35999
$ws.Range("E14").Value = 36000
$ws.Range("F14").Value = 36000

# Row 16 (Turkish)
$ws.Range("B16").Value = 4486
$ws.Range("C16").Value = 3660
$ws.Range("D16").Value = 4340
$ws.Range("E16").Value = 4476
$ws.Range("F16").Value = 4486

# Row 17 (Arabic)
$ws.Range("B17").Value = 2468
$ws.Range("C17").Value = 2432
$ws.Range("D17").Value = 2456
$ws.Range("E17").Value = 2465
$ws.Range("F17").Value = 2468

# Row 20 (Maltese)
$ws.Range("B20").Value = 595
$ws.Range("C20").Value = 427
$ws.Range("D20").Value = 580
$ws.Range("E20").Value = 592
$ws.Range("F20").Value = 594
